$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" '43.045.00'
$ws.Range("E2").Value = '  +2.56%  '
Set-TextValue "D3" '2.300.43'
$ws.Range("E3").Value = '  +1.54%  '
$ws.Range("E4").Value = '  +0.03%  '
Set-TextValue "D5" '310.19'
$ws.Range("E5").Value = '  +1.66%  '
Set-TextValue "D6" '100.53'
$ws.Range("E6").Value = '  +4.84%  '
Set-TextValue "D7" '0.536'
$ws.Range("E7").Value = '  +1.94%  '
$ws.Range("E8").Value = '  +0.02%  '
Set-TextValue "D9" '0.519'
$ws.Range("E9").Value = '  +5.97%  '
Set-TextValue "D10" '36.25'
$ws.Range("E10").Value = '  +3.32%  '
Set-TextValue "D11" '0.0825'
$ws.Range("E11").Value = '  +4.49%  '
$ws.Range("E12").Value = '  +0.79%  '
Set-TextValue "D13" '7.16'
$ws.Range("E13").Value = '  +7.87%  '
Set-TextValue "D14" '2.657.58'
$ws.Range("E14").Value = '  +1.57%  '
Set-TextValue "D15" '14.98'
$ws.Range("E15").Value = '  +4.12%  '
Set-TextValue "D16" '2.306.97'
$ws.Range("E16").Value = '  +1.81%  '
Set-TextValue "D17" '0.810'
$ws.Range("E17").Value = '  +2.13%  '
Set-TextValue "D18" '42.998.60'
$ws.Range("E18").Value = '  +2.70%  '
Set-TextValue "D19" '12.48'
$ws.Range("E19").Value = '  +1.06%  '
Set-TextValue "D20" '0.0₃0923'
$ws.Range("E20").Value = '  +2.43%  '
Set-TextValue "D21" '6.07'
$ws.Range("E21").Value = '  +1.75%  '
Set-TextValue "D22" '68.37'
$ws.Range("E22").Value = '  +0.65%  '
Set-TextValue "D23" '239.97'
$ws.Range("E23").Value = '  +1.11%  '
$ws.Range("E24").Value = '  +4.81%  '
Set-TextValue "D25" '2.62'
$ws.Range("E25").Value = '  +1.81%  '
$ws.Range("E26").Value = '  +0.06%  '
Set-TextValue "D27" '24.37'
$ws.Range("E27").Value = '  +2.99%  '
Set-TextValue "D28" '38.42'
$ws.Range("E28").Value = '  +4.50%  '
$ws.Range("B29").Value = 'Cosmos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue "D29" '9.65'
$ws.Range("E29").Value = '  +1.48%  '
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue "D30" '2.12'
$ws.Range("E30").Value = '  +0.04%  '
Set-TextValue "D31" '167.09'
$ws.Range("E31").Value = '  +4.01%  '
Set-TextValue "D32" '5.32'
$ws.Range("E32").Value = '  +2.09%  '
Set-TextValue "D33" '1.00'
$ws.Range("E33").Value = '  +0.12%  '
Set-TextValue "D34" '3.12'
$ws.Range("E34").Value = '  -1.20%  '
Set-TextValue "D35" '17.76'
$ws.Range("E35").Value = '  +3.29%  '
Set-TextValue "D36" '0.0738'
$ws.Range("E36").Value = '  +0.23%  '
$ws.Range("E37").Value = '  +0.43%  '
Set-TextValue "D38" '0.106'
$ws.Range("E38").Value = '  +1.20%  '
$ws.Range("B39").Value = 'Stellar'
$ws.Range("C39").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue "D39" '0.116'
$ws.Range("E39").Value = '  +1.62%  '
$ws.Range("B40").Value = 'ARBITRUM'
$ws.Range("C40").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue "D40" '1.83'
$ws.Range("E40").Value = '  +0.42%  '
Set-TextValue "D41" '4.21'
$ws.Range("E41").Value = '  +4.94%  '
$ws.Range("E42").Value = '  -4.42%  '
$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue "D43" '1.969.78'
$ws.Range("E43").Value = '  -0.37%  '
$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue "D44" '0.0288'
$ws.Range("E44").Value = '  +1.84%  '
Set-TextValue "D45" '19.11'
$ws.Range("E45").Value = '  +1.30%  '
Set-TextValue "D46" '3.02'
$ws.Range("E46").Value = '  +3.00%  '
Set-TextValue "D47" '9.84'
$ws.Range("E47").Value = '  -0.40%  '
Set-TextValue "D48" '55.53'
$ws.Range("E48").Value = '  +4.58%  '
$ws.Range("E49").Value = '  +15.57%  '
$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextValue "D50" '2.526.30'
$ws.Range("E50").Value = '  +1.47%  '
$ws.Range("B51").Value = 'Stacks'
$ws.Range("C51").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue "D51" '1.54'
$ws.Range("E51").Value = '  +1.84%  '
